$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.590.43"
$ws.Range("E2").Value = "  -4.57%  "
$ws.Range("D3").Value = "2.939.77"
$ws.Range("E3").Value = "  -2.34%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.98"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.44"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.511"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.45%  "
$ws.Range("D9").Value = "2.933.18"
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.126"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.39%  "
$ws.Range("E11").Value = "  -5.89%  "
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.88"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.122"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").Value = "3.424.97"
$ws.Range("E16").Value = "  -2.33%  "
$ws.Range("E17").Value = "  +5.81%  "
$ws.Range("D18").Value = "2.935.79"
$ws.Range("E18").Value = "  -2.61%  "
$ws.Range("D19").Value = "57.603.83"
$ws.Range("E19").Value = "  -4.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "417.14"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("E22").Value = "  +2.27%  "
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.10"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.80"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.46"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.46%  "
$ws.Range("E29").Value = "  +2.42%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.22"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0965"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.67"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.935"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.27"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.81%  "
$ws.Range("D38").Value = "0.0₃0687"
$ws.Range("E38").Value = "  +2.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.68"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("E40").Value = "  +3.59%  "
$ws.Range("E41").Value = "  -0.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "376.81"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0345"
$ws.Range("D43").ClearFormats()
$ws.Range("D44").Value = "2.701.11"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.38"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("E48").Value = "  +1.13%  "
$ws.Range("E49").Value = "  -2.17%  "
$ws.Range("E50").Value = "  -1.84%  "
$ws.Range("E51").Value = "  -0.34%  "
